# Fill in the worksheet with the new header / body text, replacing the old
# numeric "multiplication-table"-ish sample data with labelled string data.
#
# Cells are written in the same order the values were originally typed in
# (row 1 left-to-right, then each subsequent row's columns in the order
# A, B, C -- except C1, which reuses the text already entered in A1, and C4,
# which reuses the text already entered in B4) so the shared-string table
# comes out in the same order as the source workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "HolyFuck"
$ws.Range("B1").Value = "ShitBalls"
$ws.Range("C1").Value = "HolyFuck"

$ws.Range("A2").Value = "C1R1"
$ws.Range("A3").Value = "C1R2"
$ws.Range("A4").Value = "C1R3"

$ws.Range("B2").Value = "C2R1"
$ws.Range("B3").Value = "C2R2"
$ws.Range("B4").Value = "C3R3"

$ws.Range("C2").Value = "C3R1"
$ws.Range("C3").Value = "C3R2"
$ws.Range("C4").Value = "C3R3"

$ws.Range("A5").Value = "C1R4"
$ws.Range("B5").Value = "C2R3"
$ws.Range("C5").Value = "C3R4"

# Mark the range as having (cleared) phonetic info -- mirrors the
# noConversion phoneticPr Excel stamps on the sheet the first time phonetics
# are touched on a selection.
$null = $ws.Range("A1:C5").SetPhonetic()

# Leave the selection where the author left off, one column further right
# than the old data range.
$null = $ws.Range("E4").Select()
